$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.868.26"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "'1.887.79"
$ws.Range("E3").Value = "  -0.38%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "'0.7678"
$ws.Range("E5").Value = "  -2.18%  "

$ws.Range("D6").Value = "'242.56"
$ws.Range("E6").Value = "  -0.67%  "

$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D8").Value = "'0.3131"
$ws.Range("E8").Value = "  -0.68%  "

$ws.Range("D9").Value = "'25.59"
$ws.Range("E9").Value = "  +0.51%  "

$ws.Range("D10").Value = "'0.07120"
$ws.Range("E10").Value = "  -2.97%  "

$ws.Range("D11").Value = "'0.08534"
$ws.Range("E11").Value = "  +4.98%  "

$ws.Range("D12").Value = "'0.7628"
$ws.Range("E12").Value = "  -0.88%  "

$ws.Range("D13").Value = "'1.908.44"
$ws.Range("E13").Value = "  +3.88%  "

$ws.Range("D14").Value = "'5.364"
$ws.Range("E14").Value = "  -2.28%  "

$ws.Range("D15").Value = "'93.56"
$ws.Range("E15").Value = "  +0.10%  "

$ws.Range("D16").Value = "'6.135"
$ws.Range("E16").Value = "  -1.08%  "

$ws.Range("D17").Value = "'29.929.47"
$ws.Range("E17").Value = "  +0.34%  "

$ws.Range("D18").Value = "'13.74"
$ws.Range("E18").Value = "  -1.66%  "

$ws.Range("D19").Value = "'244.02"
$ws.Range("E19").Value = "  -0.90%  "

$ws.Range("D20").Value = "'0.000007816"

$ws.Range("D21").Value = "'0.9998"
$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("D22").Value = "'8.010"
$ws.Range("E22").Value = "  -1.66%  "

$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("D24").Value = "'0.1626"
$ws.Range("E24").Value = "  +2.32%  "

$ws.Range("D25").Value = "'9.388"
$ws.Range("E25").Value = "  -0.85%  "

$ws.Range("D26").Value = "'163.19"
$ws.Range("E26").Value = "  +0.54%  "

$ws.Range("D28").Value = "'2.035"
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("E29").Value = "  +3.44%  "

$ws.Range("D30").Value = "'1.538"
$ws.Range("E30").Value = "  -0.41%  "

$ws.Range("D31").Value = "'4.512"
$ws.Range("E31").Value = "  +0.53%  "

$ws.Range("D32").Value = "'4.127"
$ws.Range("E32").Value = "  +0.98%  "

$ws.Range("D33").Value = "'0.05441"
$ws.Range("E33").Value = "  -2.96%  "

$ws.Range("D34").Value = "'1.242"
$ws.Range("E34").Value = "  -0.98%  "

$ws.Range("D35").Value = "'0.7451"
$ws.Range("E35").Value = "  -1.49%  "

$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").Value = "'2.706"
$ws.Range("E37").Value = "  +2.28%  "

$ws.Range("D38").Value = "'0.01947"
$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("D39").Value = "'2.780"
$ws.Range("E39").Value = "  -0.19%  "

$ws.Range("D40").Value = "'0.4470"
$ws.Range("E40").Value = "  +0.07%  "

$ws.Range("D41").Value = "'1.100.94"
$ws.Range("E41").Value = "  -3.59%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.080"
$ws.Range("E42").Value = "  +1.86%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'73.04"
$ws.Range("E43").Value = "  -0.93%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = "  +0.15%  "

$ws.Range("D46").Value = "'102.98"
$ws.Range("E46").Value = "  +0.89%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'1.870"
$ws.Range("E47").Value = "  -1.65%  "

$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.665"
$ws.Range("E48").Value = "  +1.77%  "

$ws.Range("D49").Value = "'3.049"
$ws.Range("E49").Value = "  -2.48%  "

$ws.Range("D50").Value = "'2.039.05"
$ws.Range("E50").Value = "  +2.05%  "

$ws.Range("E51").Value = "  +0.34%  "
